# Tutorial 6 solution update:
#  - Column A date strings: change "/" separators to "-" separators
#    (e.g. 28/07/2022 -> 28-07-2022). Some of these (day <= 12) would be
#    mis-parsed by Excel as real dates once the separator becomes "-",
#    so for those we force Text number format first so the value is
#    kept as a literal string, matching the original "inline string"
#    storage.
#  - Attendance counters for the rows that fall on 25-08-2022 (was row
#    11), 05-09-2022 (was row 14) and 26-09-2022 (was row 20): the
#    "Total Attendance Count" (D) and "Real" (E) columns flip 0 -> 1
#    while "Absent" (H) flips 1 -> 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRows = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
$dates = @(
    "28-07-2022",
    "01-08-2022",
    "04-08-2022",
    "08-08-2022",
    "11-08-2022",
    "15-08-2022",
    "18-08-2022",
    "22-08-2022",
    "25-08-2022",
    "29-08-2022",
    "01-09-2022",
    "05-09-2022",
    "08-09-2022",
    "12-09-2022",
    "15-09-2022",
    "19-09-2022",
    "22-09-2022",
    "26-09-2022",
    "29-09-2022"
)

for ($i = 0; $i -lt $dateRows.Length; $i++) {
    $row = $dateRows[$i]
    $dateStr = $dates[$i]
    $day = [int]($dateStr.Substring(0,2))
    $cell = $ws.Cells.Item($row, 1)
    if ($day -le 12) {
        # Ambiguous as a date (day could double as a month) - keep as text.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $dateStr
}

# Attendance data corrections.
$updatedRows = @(11, 14, 20)
foreach ($row in $updatedRows) {
    $ws.Cells.Item($row, 4).Value = 1   # D: Total Attendance Count
    $ws.Cells.Item($row, 5).Value = 1   # E: Real
    $ws.Cells.Item($row, 8).Value = 0   # H: Absent
}
